# Add "Program Code" column (I) to the programs worksheet, with STC/PD/D codes
# for each program row, per commit "add modules to programs".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = "Program Code"

# Short-term Certificate rows: 2-25 -> STC100-STC123
$code = 100
for ($row = 2; $row -le 25; $row++) {
    $ws.Range("I$row").Value = "STC$code"
    $code = $code + 1
}

# Professional Diploma rows: 26-55 -> PD100-PD129
$code = 100
for ($row = 26; $row -le 55; $row++) {
    $ws.Range("I$row").Value = "PD$code"
    $code = $code + 1
}

# Degree rows: 56-62 -> D100-D106
$code = 100
for ($row = 56; $row -le 62; $row++) {
    $ws.Range("I$row").Value = "D$code"
    $code = $code + 1
}

# Adjust column widths to make room for the new column
# (Values chosen so the resulting pixel-snapped width is as close as
# possible to the authored widths of 28.48 / 20.49 / 6.51 characters.)
$ws.Columns.Item(1).ColumnWidth = 27.666666666666664
$ws.Columns.Item(2).ColumnWidth = 19.666666666666664
$ws.Columns.Item(3).ColumnWidth = 5.666666666666666

# Restore the view/selection that results from scrolling/selecting during the edit
$ws.Application.ActiveWindow.ScrollRow = 41
$ws.Range("G56").Select()
